# Adds new patient-log rows (project columns) and categorizes the MOCA-score
# note that used to live in the old row 44.
#
# Before: row 44 held ID 07122 + a long wrapped "Error may have occured..."
# note in column D, and row 45 was essentially blank.
# After:  row 44 is repurposed for ID 07129, and four additional rows
# (45-48) are appended, with the trailing note moved down to row 48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 44: replace the old ID/notes pair with the new entry -------------
# Drop the old Notes cell (D44) entirely - the new row 44 has no Notes.
$ws.Range("D44").Clear()

$ws.Range("A44").Value2 = "07129"
$ws.Range("B44").Value2 = "01/22/2020"
$ws.Range("C44").NumberFormat = "h:mm"
$ws.Range("C44").Value2 = 0.4548611111111111
$ws.Rows.Item(44).RowHeight = 15.75

# --- Row 45 -----------------------------------------------------------------
$ws.Range("A45").Value2 = "07132"
$ws.Range("B45").Value2 = "01/23/2020"
$ws.Range("C45").NumberFormat = "h:mm"
$ws.Range("C45").Value2 = 0.55069444444444449

# --- Row 46 -----------------------------------------------------------------
$ws.Range("A46").Value2 = "07135"
$ws.Range("B46").Value2 = "01/29/2020"
$ws.Range("C46").NumberFormat = "h:mm"
$ws.Range("C46").Value2 = 0.47222222222222227

# --- Row 47 -----------------------------------------------------------------
$ws.Range("A47").Value2 = "07107"
$ws.Range("B47").Value2 = "01/07/2020"
$ws.Range("C47").NumberFormat = "h:mm"
$ws.Range("C47").Value2 = 0.60763888888888895

# --- Row 48 (new Notes entry, moved down from the old row 44) ---------------
$ws.Range("A48").Value2 = "07136"
$ws.Range("B48").Value2 = "01/28/2020"
$ws.Range("C48").NumberFormat = "h:mm"
$ws.Range("C48").Value2 = 0.49444444444444446
$ws.Range("D48").Value2 = "File split in two parts. "
$ws.Rows.Item(48).RowHeight = 15.75

# --- Keep the selection where Excel would leave it after entering this data -
$ws.Range("D49").Select() | Out-Null
